# Update cryptos list: apply price (D) and volume % change (E) updates
# For numeric-looking price strings, a leading apostrophe is used so Excel
# stores them as text (preserving exact formatting like "564.92" or "1.00")
# instead of auto-converting to a floating point number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.366.88"
$ws.Range("E2").Value = "  +3.96%  "

$ws.Range("D3").Value = "2.998.05"
$ws.Range("E3").Value = "  +3.71%  "

$ws.Range("E4").Value = "  +0.10%  "

$ws.Range("D5").Value = "'564.92"
$ws.Range("E5").Value = "  +3.31%  "

$ws.Range("D6").Value = "'139.47"
$ws.Range("E6").Value = "  +13.48%  "

$ws.Range("E7").Value = "  -0.22%  "

$ws.Range("E8").Value = "  +4.92%  "

$ws.Range("D9").Value = "2.992.09"
$ws.Range("E9").Value = "  +3.81%  "

$ws.Range("E10").Value = "  +10.26%  "

$ws.Range("D11").Value = "'5.05"
$ws.Range("E11").Value = "  +9.32%  "

$ws.Range("D12").Value = "'0.461"
$ws.Range("E12").Value = "  +5.48%  "

$ws.Range("E13").Value = "  +10.08%  "

$ws.Range("D14").Value = "'33.92"
$ws.Range("E14").Value = "  +4.60%  "

$ws.Range("E15").Value = "  +2.79%  "

$ws.Range("D16").Value = "3.497.50"
$ws.Range("E16").Value = "  +3.91%  "

$ws.Range("D17").Value = "'7.02"
$ws.Range("E17").Value = "  +7.27%  "

$ws.Range("D18").Value = "2.997.84"
$ws.Range("E18").Value = "  +3.84%  "

$ws.Range("D19").Value = "59.315.97"
$ws.Range("E19").Value = "  +3.78%  "

$ws.Range("D20").Value = "'430.49"
$ws.Range("E20").Value = "  +6.88%  "

$ws.Range("D21").Value = "'13.63"
$ws.Range("E21").Value = "  +6.19%  "

$ws.Range("E22").Value = "  +6.06%  "

$ws.Range("E23").Value = "  +5.77%  "

$ws.Range("D24").Value = "'13.47"
$ws.Range("E24").Value = "  +5.66%  "

$ws.Range("D25").Value = "'80.76"
$ws.Range("E25").Value = "  +4.73%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = "  +0.33%  "

$ws.Range("E28").Value = "  +11.59%  "

$ws.Range("E29").Value = "  +4.43%  "

$ws.Range("E30").Value = "  +8.90%  "

$ws.Range("D31").Value = "'25.77"
$ws.Range("E31").Value = "  +4.66%  "

$ws.Range("D32").Value = "'6.17"
$ws.Range("E32").Value = "  +3.53%  "

$ws.Range("D33").Value = "'0.0997"
$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").Value = "'1.01"
$ws.Range("E34").Value = "  +10.69%  "

$ws.Range("D35").Value = "0.0₃0770"
$ws.Range("E35").Value = "  +24.03%  "

$ws.Range("E36").Value = "  +7.27%  "

$ws.Range("D37").Value = "'2.08"
$ws.Range("E37").Value = "  +4.50%  "

$ws.Range("D38").Value = "'49.43"
$ws.Range("E38").Value = "  +3.47%  "

$ws.Range("D39").Value = "'8.69"
$ws.Range("E39").Value = "  +5.18%  "

$ws.Range("E40").Value = "  +15.67%  "

$ws.Range("D41").Value = "'407.56"
$ws.Range("E41").Value = "  +13.77%  "

$ws.Range("E42").Value = "  +4.12%  "

$ws.Range("D43").Value = "2.762.73"
$ws.Range("E43").Value = "  +5.20%  "

$ws.Range("E44").Value = "  +3.22%  "

$ws.Range("D45").Value = "'0.251"
$ws.Range("E45").Value = "  +10.00%  "

$ws.Range("D46").Value = "'0.999"
$ws.Range("E46").Value = "  +0.03%  "

$ws.Range("D47").Value = "'124.49"
$ws.Range("E47").Value = "  +4.28%  "

$ws.Range("D48").Value = "'2.02"
$ws.Range("E48").Value = "  +4.67%  "

$ws.Range("E49").Value = "  +3.22%  "

$ws.Range("D50").Value = "'32.67"
$ws.Range("E50").Value = "  +21.06%  "

$ws.Range("D51").Value = "'23.60"
$ws.Range("E51").Value = "  +3.75%  "

